# Auto-generated PowerShell COM-interop script to apply the cryptos.xlsx update
# Commit: Updated cryptos list on Sun Aug 11 22:57:32 UTC 2024 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.007.27'
$ws.Range('E2').Value = '  -3.03%  '
$ws.Range('D3').Value = '2.565.93'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '507.51'
$ws.Range('E5').Value = '  -3.10%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.00'
$ws.Range('E6').Value = '  -7.65%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.554'
$ws.Range('E8').Value = '  -5.71%  '
$ws.Range('D9').Value = '2.571.23'
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.22'
$ws.Range('E10').Value = '  -7.11%  '
$ws.Range('E11').Value = '  -3.64%  '
$ws.Range('E12').Value = '  -4.78%  '
$ws.Range('E13').Value = '  -0.97%  '
$ws.Range('D14').Value = '3.009.72'
$ws.Range('E14').Value = '  -1.57%  '
$ws.Range('D15').Value = '58.963.39'
$ws.Range('E15').Value = '  -3.11%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '20.60'
$ws.Range('E16').Value = '  -4.80%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000135'
$ws.Range('E17').Value = '  -4.83%  '
$ws.Range('D18').Value = '2.564.96'
$ws.Range('E18').Value = '  -1.63%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.52'
$ws.Range('E19').Value = '  -5.02%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '332.15'
$ws.Range('E20').Value = '  -6.53%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.07'
$ws.Range('E21').Value = '  -4.78%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.95'
$ws.Range('E23').Value = '  -4.08%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '59.50'
$ws.Range('E24').Value = '  -2.41%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.406'
$ws.Range('E25').Value = '  -4.69%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.157'
$ws.Range('E27').Value = '  -5.78%  '
$ws.Range('D28').Value = '0.0₃0777'
$ws.Range('E28').Value = '  -8.12%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.87'
$ws.Range('E29').Value = '  -7.09%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '149.90'
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '18.59'
$ws.Range('E32').Value = '  -4.22%  '
$ws.Range('E34').Value = '  -3.79%  '
$ws.Range('E35').Value = '  -6.96%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.890'
$ws.Range('E36').Value = '  -2.89%  '
$ws.Range('E37').Value = '  -7.99%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '35.92'
$ws.Range('E38').Value = '  -1.53%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.827'
$ws.Range('E39').Value = '  -9.04%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '287.41'
$ws.Range('E40').Value = '  -1.44%  '
$ws.Range('E41').Value = '  -8.10%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.50'
$ws.Range('E42').Value = '  -8.15%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.608'
$ws.Range('E44').Value = '  -2.28%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0980'
$ws.Range('E45').Value = '  -3.08%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0531'
$ws.Range('E46').Value = '  -5.11%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.35'
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('E48').Value = '  -5.00%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0227'
$ws.Range('E49').Value = '  -4.62%  '
$ws.Range('E50').Value = '  -8.04%  '
$ws.Range('D51').Value = '1.912.33'
$ws.Range('E51').Value = '  -2.33%  '
